$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15, pushing existing rows 15..120 down to 16..121
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15 with its values
$ws.Range("A15").Value = 4
$ws.Range("B15").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C15").Value = "Los Lagos"
$ws.Range("D15").Value = 44462
$ws.Range("E15").Value = 10
$ws.Range("F15").Value = 100112032
$ws.Range("G15").Value = "Zapallo italiano"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 17000
$ws.Range("L15").Value = 17000
$ws.Range("M15").Value = 17000
$ws.Range("N15").Value = "$/caja 50 unidades"
$ws.Range("O15").Value = "Región de Arica y Parinacota"
$ws.Range("P15").Value = 340
$ws.Range("Q15").Value = 50
$ws.Range("R15").Value = "Hortaliza"
